# Actualizacion automatica del mapa (2025-11-13 12:03:39)
# The source row for "Caso" 7689 (VIALE, LUIS 2026) is removed; all rows below it shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 79 corresponds to Caso 7689 / VIALE, LUIS 2026 - delete it entirely, shifting rows below up.
$ws.Rows.Item(79).Delete()
